# Apply the "Adding Master Data XLS" edit:
#  - Replace the Dongle (DNG) rows with Desktop Computer (DKS) rows
#  - Update the active cell selection on the sheet
#  - Set page setup (paper size / orientation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 : DKS / Dekstop / Desktop Computer / eng
$ws.Range("A2").Value = "DKS"
$ws.Range("B2").Value = "Dekstop"
$ws.Range("C2").Value = "Desktop Computer"
$ws.Range("D2").Value = "eng"

# Row 3 : DKS / الحاسوب / أجهزة الكمبيوتر المكتبية / ara
$ws.Range("A3").Value = "DKS"
$ws.Range("B3").Value = "الحاسوب"
$ws.Range("C3").Value = "أجهزة الكمبيوتر المكتبية"
$ws.Range("D3").Value = "ara"

# Row 4 : DKS / Ordinateur / Ordinateurs de bureau / fra
$ws.Range("A4").Value = "DKS"
$ws.Range("B4").Value = "Ordinateur"
$ws.Range("C4").Value = "Ordinateurs de bureau"
$ws.Range("D4").Value = "fra"

# Update the saved selection on the worksheet
$ws.Range("D10").Select()

# Set the page setup (paper size = A4(9), orientation = portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
